$wb = $excel.ActiveWorkbook

function Set-LinkDisplay {
    param($ws, $addr, $newDisplay)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $newDisplay
        }
    }
}

# ---------------------------------------------------------------
# Sheet "Overview": the ccd37122... file now sorts to row 2, the
# 24fb51b8... file moves to row 3. Row 3 status/date now reflect a
# fresh "Ready for handoff" report generation.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsOverview.Range("B2").Value = "e2e\ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsOverview.Range("A3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
$wsOverview.Range("B3").Value = "e2e\24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 22:47:27"

Set-LinkDisplay $wsOverview '$B$2' "e2e\ccd37122-6c86-4e63-b00d-1275dcf94100.md"
Set-LinkDisplay $wsOverview '$B$3' "e2e\24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"

# ---------------------------------------------------------------
# Sheet "zh-cn": same row re-identification; row 2 = ccd, row 3 =
# 24fb. Row 3 (24fb) picks up a new handoff datetime and an error
# about the handback file being stale.
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("G2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.0ce5dcc0d556b017a9782c5d3538e6875f3f81f3.zh-cn.xlf"
$wsZh.Range("I2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsZh.Range("J2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.0ce5dcc0d556b017a9782c5d3538e6875f3f81f3.zh-cn.xlf"

$wsZh.Range("A3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.ba29efadbd1379a9a4fe3495af58894cb7218b4c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-22 22:47:22"
$wsZh.Range("I3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
$wsZh.Range("J3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.ba29efadbd1379a9a4fe3495af58894cb7218b4c.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c0e7c8d7ecdb1df965a13ab9224c2dc0bc62175/e2e/24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87f5590e9a6fd269bf101c9cf4b49580cbb822cc/e2e/24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md."

Set-LinkDisplay $wsZh '$A$2' "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
Set-LinkDisplay $wsZh '$I$2' "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
Set-LinkDisplay $wsZh '$A$3' "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
Set-LinkDisplay $wsZh '$I$3' "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------
# Sheet "de-de": same row re-identification as zh-cn.
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsDe.Range("G2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.0ce5dcc0d556b017a9782c5d3538e6875f3f81f3.de-de.xlf"
$wsDe.Range("I2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
$wsDe.Range("J2").Value = "ccd37122-6c86-4e63-b00d-1275dcf94100.0ce5dcc0d556b017a9782c5d3538e6875f3f81f3.de-de.xlf"

$wsDe.Range("A3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.ba29efadbd1379a9a4fe3495af58894cb7218b4c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-22 22:47:27"
$wsDe.Range("I3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
$wsDe.Range("J3").Value = "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.ba29efadbd1379a9a4fe3495af58894cb7218b4c.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c0e7c8d7ecdb1df965a13ab9224c2dc0bc62175/e2e/24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87f5590e9a6fd269bf101c9cf4b49580cbb822cc/e2e/24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md."

Set-LinkDisplay $wsDe '$A$2' "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
Set-LinkDisplay $wsDe '$I$2' "ccd37122-6c86-4e63-b00d-1275dcf94100.md"
Set-LinkDisplay $wsDe '$A$3' "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"
Set-LinkDisplay $wsDe '$I$3' "24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md"

$wsDe.Columns.Item(16).ColumnWidth = 39.17
